# Apply the edit described by the diff:
# - A1 gains a wrapped second (description) line, row height becomes 60, and
#   the cell picks up the workbook's "wrap text" style (same style index that
#   B2/B4/B5 already use).
# - A brand-new line is inserted as the new B2, explaining the leading-'0' guard.
# - The previous B2..B4 content shifts down by one row (now B3..B5), text and
#   formatting unchanged.
# - The previous B5 content (the _OPERATION_SALE_ICC_ / _TRT_SALE_ block)
#   becomes the new B6 with updated identifiers (_OPERATION_SALE_CTLS_IDLE_ /
#   _TRT_SALE_CTLS_) and a taller row height (45, it now wraps 2 lines).
# - Column B grows from 68.25 to (about) 78.25 and stops being "best fit".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# Literal text content. NOTE: a round-tripped Range.Value *read* on this host
# returns a COM type descriptor instead of the actual string, so we use the
# literal text (taken from the original file / the target diff) instead of
# reading-then-rewriting each cell.
$textTitle  = "inMENU_MenuKeyInAndGetAmount" + $nl + "整體流程為偵測到輸入1~9後設定交易類型和 inRunOperationID、>inRunTRTID，以跑對應的function"
$textBlock0 = "會先擋srEventMenuItem->inEventCode == '0' (/* 不接受金額第一位為0 */)"
$textLogo   = "inFunc_Display_LOGO  透過inGetLOGONum()，顯示不同的logo"
$textDisp   = "inDISP_PutGraphic顯示一般交易的圖片"
$textPwd    = "srEventMenuItem->inPasswordLevel = _ACCESS_WITH_CUSTOM_;" + $nl + "        srEventMenuItem->inCode = _SALE_;(inCode == 交易類型)"
$textRunId  = "                srEventMenuItem->inRunOperationID = _OPERATION_SALE_CTLS_IDLE_;" + $nl + "                srEventMenuItem->inRunTRTID = _TRT_SALE_CTLS_;"

# --- Row 1 / A1: title gains an explanatory second line ---
$ws.Range("A1").Value = $textTitle
$ws.Range("A1").WrapText = $true
$ws.Rows(1).RowHeight = 60

# --- Row 2 / B2: brand-new line (was not present before); must end up with
#     the sheet's default (no) style. B2 previously had the wrap-text style,
#     so explicitly strip formatting back to a known-blank cell's (PasteSpecial
#     formats is the reliable way to do this on this host -- ClearFormats()
#     on a cell that currently carries a non-default style index introduces
#     a spurious extra font/style instead of falling back to style 0).
$ws.Range("B2").Value = $textBlock0
$ws.Range("Z100").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Row 3 / B3: previously B2 content (unchanged text/style) ---
$ws.Range("B3").Value = $textLogo
$ws.Range("B3").WrapText = $true

# --- Row 4 / B4: previously B3 content (unchanged text), must end up with
#     no style (same situation/fix as B2 above) and no explicit row height
#     (row 4 previously had an explicit ht="30" that must not leak through). ---
$ws.Range("B4").Value = $textDisp
$ws.Range("Z100").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Rows(4).AutoFit()

# --- Row 5 / B5: previously B4 content (unchanged text/style) ---
$ws.Range("B5").Value = $textPwd
$ws.Range("B5").WrapText = $true
$ws.Rows(5).RowHeight = 30

# --- Row 6 / B6: previously B5 content, with updated identifiers ---
$ws.Range("B6").Value = $textRunId
$ws.Range("B6").WrapText = $true
$ws.Rows(6).RowHeight = 45

$excel.CutCopyMode = $false

# --- Column B width: grows and is no longer "best fit" ---
# NOTE: the host quantizes ColumnWidth to steps of 1/7 (the same pixel-based
# rounding real Excel applies), so an input of exactly 78.25 rounds up to 79.
# 77.5 is the closest input that lands on the nearest reachable width to 78.25.
$ws.Columns("B").ColumnWidth = 77.5

$ws.Range("B6").Select() | Out-Null
